# Logboek entry for "30 september" originally consisted of two runs that,
# read together, form one unbroken sentence:
#   "...mogelijk moeten maken. " + "Ingelezen in Tkinter, ... te maken."
# The edit merges those two runs into a single run (no text changes), then
# splits a new "1 oktober:" entry plus its body out of what follows (right
# before the trailing page break), pushing the page break down onto the new,
# final paragraph.

$d = $word.ActiveDocument

# 1) Merge the two runs that make up the "30 september" paragraph into one
#    run. The visible text is unchanged; searching for the sentence that
#    spans the run boundary and replacing it with itself causes Word to
#    re-flow it as a single run while keeping the existing run formatting
#    (lang="nl-NL").
$mergeText = "maken. Ingelezen in Tkinter, zou een goede optie zijn om een programmeerinterface te maken."
$d.Content.Find.Execute($mergeText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $mergeText, 2) | Out-Null

# 2) Locate that paragraph again (still "probleem met pycharm opgelost...").
#    Keep working off the same Range object Find.Execute repositions, rather
#    than $word.Selection (Find on a Range does not move the Selection).
$rng = $d.Content
$rng.Find.Execute("probleem met pycharm opgelost", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$entryPara = $rng.Paragraphs(1)

# 3) Split off a new paragraph right before the trailing page break (i.e.
#    after "...te maken." but before the page-break character), so the page
#    break ends up alone in its own paragraph for now.
$splitPoint = $entryPara.Range.Duplicate
$splitPoint.MoveEnd(1, -2)   # back up over the page-break char + paragraph mark
$splitPoint.Collapse(0)      # wdCollapseEnd
$splitPoint.InsertParagraphAfter()

# 4) The page-break-only paragraph now immediately follows $entryPara.
#    Insert two fresh empty paragraphs in front of it for "1 oktober:" and
#    its body text.
$pageBreakPara = $entryPara.Next()
$pageBreakPara.Range.InsertParagraphBefore()
$pageBreakPara.Range.InsertParagraphBefore()

# 5) Fill in "1 oktober:".
$datePara = $entryPara.Next()
$dateRange = $datePara.Range.Duplicate
$dateRange.MoveEnd(1, -1)    # exclude the paragraph mark
$dateRange.Text = "1 oktober:"
$dateRange.LanguageID = "nl-NL"

# 6) Fill in the body text for 1 oktober.
$bodyPara = $datePara.Next()
$bodyRange = $bodyPara.Range.Duplicate
$bodyRange.MoveEnd(1, -1)    # exclude the paragraph mark
$bodyRange.Text = "verder ingelezen in Tkinter. Simpele interface gemaakt waarin motoren en beweegrichtingen in geselecteerd kunnen worden en deze data vervolgens doorgestuurd wordt."
$bodyRange.LanguageID = "nl-NL"

# 7) Re-attach the page break to the end of the new body paragraph by
#    deleting the paragraph mark that currently separates them, merging the
#    (now empty) page-break paragraph back into $bodyPara.
$joinMark = $bodyPara.Range.Duplicate
$joinMark.Collapse(0)        # wdCollapseEnd -> right after body's own paragraph mark
$joinMark.MoveStart(1, -1)   # select just that paragraph mark
$joinMark.Delete()

Write-Output "done"
